$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column N with year 2022 data, mirroring column M's styling.

# N2: empty cell, same style as M2 (s="3")
$ws.Range("N2").Value = $null
$ws.Range("N2").Style = $ws.Range("M2").Style

# N3: header year 2022, same style as M3 (s="13")
$ws.Range("N3").Value = 2022
$ws.Range("N3").Style = $ws.Range("M3").Style

# N4: value 1434, same style as M4 (s="14")
$ws.Range("N4").Value = 1434
$ws.Range("N4").Style = $ws.Range("M4").Style

# N5: value 12822, same style as M5 (s="14")
$ws.Range("N5").Value = 12822
$ws.Range("N5").Style = $ws.Range("M5").Style

# N6: value 3099, same style as M6 (s="14")
$ws.Range("N6").Value = 3099
$ws.Range("N6").Style = $ws.Range("M6").Style

# N7: value 9722, same style as M7 (s="14")
$ws.Range("N7").Value = 9722
$ws.Range("N7").Style = $ws.Range("M7").Style

# N8: value 14424, same style as M8 (s="15")
$ws.Range("N8").Value = 14424
$ws.Range("N8").Style = $ws.Range("M8").Style

# N9: value 5279, same style as M9 (s="15")
$ws.Range("N9").Value = 5279
$ws.Range("N9").Style = $ws.Range("M9").Style

# N10: value 9145, same style as M10 (s="16")
$ws.Range("N10").Value = 9145
$ws.Range("N10").Style = $ws.Range("M10").Style

# Update the selection to N2 (mirrors the author's last-click position)
$ws.Range("N2").Select()
